# Miami Heat roster update:
#  - Gabe Vincent and Jimmy Butler move up in roster order (ahead of Kyle Lowry)
#  - Kyle Lowry shifts down to where Jimmy Butler used to be
#  - Jamaree Bouyea leaves the roster; Omer Yurtseven's row data shifts up into his slot
#  - Kevin Love joins the roster, taking the final row (where Omer Yurtseven used to be)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: becomes Gabe Vincent
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "Gabe Vincent"
$ws.Range("D6").Value = "PG"
$ws.Range("E6").Value = "6-3"
$ws.Range("F6").Value = 200
$ws.Range("G6").Value = "June 14, 1996"
$ws.Range("H6").Value = "us"
$ws.Range("I6").Value = "3"
$ws.Range("J6").Value = "UC Santa Barbara"
$ws.Range("K6").Value = "https://www.basketball-reference.com/players/v/vincega01.html"

# Row 7: becomes Jimmy Butler
$ws.Range("B7").Value = 22
$ws.Range("C7").Value = "Jimmy Butler"
$ws.Range("D7").Value = "SF"
$ws.Range("E7").Value = "6-7"
$ws.Range("F7").Value = 230
$ws.Range("G7").Value = "September 14, 1989"
$ws.Range("H7").Value = "us"
$ws.Range("I7").Value = "11"
$ws.Range("J7").Value = "Marquette"
$ws.Range("K7").Value = "https://www.basketball-reference.com/players/b/butleji01.html"

# Row 8: becomes Kyle Lowry
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = "Kyle Lowry"
$ws.Range("D8").Value = "PG"
$ws.Range("E8").Value = "6-0"
$ws.Range("F8").Value = 196
$ws.Range("G8").Value = "March 25, 1986"
$ws.Range("H8").Value = "us"
$ws.Range("I8").Value = "16"
$ws.Range("J8").Value = "Villanova"
$ws.Range("K8").Value = "https://www.basketball-reference.com/players/l/lowryky01.html"

# Row 16: Jamaree Bouyea leaves, replaced by Omer Yurtseven's data (no jersey number)
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "Omer Yurtseven"
$ws.Range("D16").Value = "C"
$ws.Range("E16").Value = "7-0"
$ws.Range("F16").Value = 264
$ws.Range("G16").Value = "June 19, 1998"
$ws.Range("H16").Value = "tr"
$ws.Range("I16").Value = "1"
$ws.Range("J16").Value = "NC State, Georgetown"
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/y/yurtsom01.html"

# Row 17: new player Kevin Love joins (no jersey number)
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "Kevin Love"
$ws.Range("D17").Value = "PF"
$ws.Range("E17").Value = "6-8"
$ws.Range("F17").Value = 251
$ws.Range("G17").Value = "September 7, 1988"
$ws.Range("H17").Value = "us"
$ws.Range("I17").Value = "14"
$ws.Range("J17").Value = "UCLA"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/l/loveke01.html"
